$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fix the dialogue text in A2: "I like the most quiet restaurants."
#        -> "I prefer quiet restaurants." and drop the trailing newline
#        that used to force xml:space="preserve" on the shared string.
$old = $ws.Range("A2").Value2
$new = $old.Replace("I like the most quiet restaurants.", "I prefer quiet restaurants.")
if ($new.EndsWith("`n")) {
    $new = $new.Substring(0, $new.Length - 1)
}
$ws.Range("A2").Value = $new

# --- 2. Scroll the first sheet's view so A2 is the top-left visible cell
#        (mirrors <sheetView tabSelected="1" topLeftCell="A2" .../> in the xml).
$ws.Activate()
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
